$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.001.05"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "3.515.70"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'591.12"
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("D6").Value = "'133.78"
$ws.Range("E6").Value = "  -1.27%  "
$ws.Range("D7").Value = "3.514.68"
$ws.Range("E7").Value = "  -1.43%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("E10").Value = "  +1.31%  "
$ws.Range("D11").Value = "'7.20"
$ws.Range("E11").Value = "  +3.11%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "4.115.26"
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("D14").Value = "'27.66"
$ws.Range("E14").Value = "  +2.27%  "
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "3.515.78"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").Value = "65.001.01"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").Value = "'14.31"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("E21").Value = "  -2.77%  "
$ws.Range("D22").Value = "'390.97"
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("D23").Value = "'0.577"
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("D24").Value = "3.659.05"
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("D25").Value = "'74.68"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -3.43%  "
$ws.Range("E28").Value = "  +9.47%  "
$ws.Range("D29").Value = "'7.63"
$ws.Range("E29").Value = "  -1.27%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("D32").Value = "'8.30"
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("D33").Value = "3.523.68"
$ws.Range("E33").Value = "  -1.27%  "
$ws.Range("D34").Value = "'24.05"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("D37").Value = "'5.23"
$ws.Range("E37").Value = "  +4.92%  "
$ws.Range("E38").Value = "  +1.49%  "
$ws.Range("D39").Value = "'169.47"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "'6.91"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("D43").Value = "'26.02"
$ws.Range("E43").Value = "  -4.23%  "
$ws.Range("E44").Value = "  +3.61%  "
$ws.Range("D45").Value = "'42.92"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.445.58"
$ws.Range("E49").Value = "  -1.41%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'6.87"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").Value = "'0.900"
$ws.Range("E51").Value = "  +3.26%  "
